$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: D2: 44510 -> 44453, J2: 40 -> 50, K2: 15000 -> 14000, L2: 16000 -> 15000, M2: 15500 -> 14600, P2: 1192 -> 1123
$ws.Range("D2").Value = 44453
$ws.Range("J2").Value = 50
$ws.Range("K2").Value = 14000
$ws.Range("L2").Value = 15000
$ws.Range("M2").Value = 14600
$ws.Range("P2").Value = 1123

# Row 3: D3: 44813 -> 44425, J3: 50 -> 60, K3: 13000 -> 14000, L3: 14000 -> 15000, M3: 13400 -> 14500, P3: 1031 -> 1115
$ws.Range("D3").Value = 44425
$ws.Range("J3").Value = 60
$ws.Range("K3").Value = 14000
$ws.Range("L3").Value = 15000
$ws.Range("M3").Value = 14500
$ws.Range("P3").Value = 1115

# Row 4: D4: 44435 -> 44610, J4: 100 -> 50, K4: 13000 -> 17000, L4: 14000 -> 18000, M4: 13500 -> 17400, P4: 1038 -> 1338
$ws.Range("D4").Value = 44610
$ws.Range("J4").Value = 50
$ws.Range("K4").Value = 17000
$ws.Range("L4").Value = 18000
$ws.Range("M4").Value = 17400
$ws.Range("P4").Value = 1338

# Row 6: D6: 44308 -> 44334, L6: 27000 -> 28000, M6: 26400 -> 27200, P6: 2031 -> 2092
$ws.Range("D6").Value = 44334
$ws.Range("L6").Value = 28000
$ws.Range("M6").Value = 27200
$ws.Range("P6").Value = 2092

# Row 7: D7: 44453 -> 44350, J7: 50 -> 40, K7: 14000 -> 23000, L7: 15000 -> 25000, M7: 14600 -> 24000, P7: 1123 -> 1846
$ws.Range("D7").Value = 44350
$ws.Range("J7").Value = 40
$ws.Range("K7").Value = 23000
$ws.Range("L7").Value = 25000
$ws.Range("M7").Value = 24000
$ws.Range("P7").Value = 1846

# Row 8: D8: 44320 -> 44474, J8: 50 -> 40, K8: 26000 -> 13000, L8: 28000 -> 14000, M8: 26800 -> 13500, P8: 2062 -> 1038
$ws.Range("D8").Value = 44474
$ws.Range("J8").Value = 40
$ws.Range("K8").Value = 13000
$ws.Range("L8").Value = 14000
$ws.Range("M8").Value = 13500
$ws.Range("P8").Value = 1038

# Row 9: D9: 44433 -> 44769, J9: 100 -> 50, K9: 13000 -> 14000, L9: 14000 -> 15000, M9: 13500 -> 14600, P9: 1038 -> 1123
$ws.Range("D9").Value = 44769
$ws.Range("J9").Value = 50
$ws.Range("K9").Value = 14000
$ws.Range("L9").Value = 15000
$ws.Range("M9").Value = 14600
$ws.Range("P9").Value = 1123

# Row 11: D11: 44708 -> 44813, M11: 13600 -> 13400, P11: 1046 -> 1031
$ws.Range("D11").Value = 44813
$ws.Range("M11").Value = 13400
$ws.Range("P11").Value = 1031

# Row 12: D12: 44719 -> 44691, J12: 50 -> 100, K12: 13000 -> 12000, L12: 14000 -> 13000, M12: 13400 -> 12500, P12: 1031 -> 962
$ws.Range("D12").Value = 44691
$ws.Range("J12").Value = 100
$ws.Range("K12").Value = 12000
$ws.Range("L12").Value = 13000
$ws.Range("M12").Value = 12500
$ws.Range("P12").Value = 962

# Row 13: D13: 44350 -> 44362, K13: 23000 -> 15000, L13: 25000 -> 16000, M13: 24000 -> 15500, P13: 1846 -> 1192
$ws.Range("D13").Value = 44362
$ws.Range("K13").Value = 15000
$ws.Range("L13").Value = 16000
$ws.Range("M13").Value = 15500
$ws.Range("P13").Value = 1192

# Row 14: D14: 44610 -> 44664, K14: 17000 -> 11000, L14: 18000 -> 12000, M14: 17400 -> 11600, P14: 1338 -> 892
$ws.Range("D14").Value = 44664
$ws.Range("K14").Value = 11000
$ws.Range("L14").Value = 12000
$ws.Range("M14").Value = 11600
$ws.Range("P14").Value = 892

# Row 15: D15: 44664 -> 44761, J15: 50 -> 25, K15: 11000 -> 14000, L15: 12000 -> 15000, M15: 11600 -> 14400, P15: 892 -> 1108
$ws.Range("D15").Value = 44761
$ws.Range("J15").Value = 25
$ws.Range("K15").Value = 14000
$ws.Range("L15").Value = 15000
$ws.Range("M15").Value = 14400
$ws.Range("P15").Value = 1108

# Row 16: D16: 44355 -> 44159, K16: 18000 -> 30000, L16: 20000 -> 32000, M16: 19000 -> 31000, P16: 1462 -> 2385
$ws.Range("D16").Value = 44159
$ws.Range("K16").Value = 30000
$ws.Range("L16").Value = 32000
$ws.Range("M16").Value = 31000
$ws.Range("P16").Value = 2385

# Row 17: D17: 44761 -> 44435, J17: 25 -> 100, K17: 14000 -> 13000, L17: 15000 -> 14000, M17: 14400 -> 13500, P17: 1108 -> 1038
$ws.Range("D17").Value = 44435
$ws.Range("J17").Value = 100
$ws.Range("K17").Value = 13000
$ws.Range("L17").Value = 14000
$ws.Range("M17").Value = 13500
$ws.Range("P17").Value = 1038

# Row 18: D18: 44362 -> 44509, J18: 40 -> 100
$ws.Range("D18").Value = 44509
$ws.Range("J18").Value = 100

# Row 19: D19: 44705 -> 44316, K19: 10000 -> 27000, L19: 11000 -> 28000, M19: 10400 -> 27400, P19: 800 -> 2108
$ws.Range("D19").Value = 44316
$ws.Range("K19").Value = 27000
$ws.Range("L19").Value = 28000
$ws.Range("M19").Value = 27400
$ws.Range("P19").Value = 2108

# Row 20: D20: 44334 -> 44775, J20: 50 -> 20, K20: 26000 -> 12000, L20: 28000 -> 13000, M20: 27200 -> 12500, P20: 2092 -> 962
$ws.Range("D20").Value = 44775
$ws.Range("J20").Value = 20
$ws.Range("K20").Value = 12000
$ws.Range("L20").Value = 13000
$ws.Range("M20").Value = 12500
$ws.Range("P20").Value = 962

# Row 21: D21: 44377 -> 44503, J21: 40 -> 35, K21: 14000 -> 15000, L21: 15000 -> 16000, M21: 14500 -> 15429, P21: 1115 -> 1187
$ws.Range("D21").Value = 44503
$ws.Range("J21").Value = 35
$ws.Range("K21").Value = 15000
$ws.Range("L21").Value = 16000
$ws.Range("M21").Value = 15429
$ws.Range("P21").Value = 1187

# Row 22: D22: 44313 -> 44523, J22: 50 -> 40, K22: 25000 -> 15000, L22: 26000 -> 16000, M22: 25600 -> 15500, P22: 1969 -> 1192
$ws.Range("D22").Value = 44523
$ws.Range("J22").Value = 40
$ws.Range("K22").Value = 15000
$ws.Range("L22").Value = 16000
$ws.Range("M22").Value = 15500
$ws.Range("P22").Value = 1192

# Row 23: D23: 44383 -> 44708, K23: 15000 -> 13000, L23: 16000 -> 14000, M23: 15400 -> 13600, P23: 1185 -> 1046
$ws.Range("D23").Value = 44708
$ws.Range("K23").Value = 13000
$ws.Range("L23").Value = 14000
$ws.Range("M23").Value = 13600
$ws.Range("P23").Value = 1046

# Row 24: D24: 44509 -> 44741, J24: 100 -> 50, K24: 15000 -> 14000, L24: 16000 -> 15000, M24: 15500 -> 14400, P24: 1192 -> 1108
$ws.Range("D24").Value = 44741
$ws.Range("J24").Value = 50
$ws.Range("K24").Value = 14000
$ws.Range("L24").Value = 15000
$ws.Range("M24").Value = 14400
$ws.Range("P24").Value = 1108

# Row 25: D25: 44503 -> 44755, J25: 35 -> 40, K25: 15000 -> 14000, L25: 16000 -> 15000, M25: 15429 -> 14500, P25: 1187 -> 1115
$ws.Range("D25").Value = 44755
$ws.Range("J25").Value = 40
$ws.Range("K25").Value = 14000
$ws.Range("L25").Value = 15000
$ws.Range("M25").Value = 14500
$ws.Range("P25").Value = 1115

# Row 26: D26: 44264 -> 44313, J26: 40 -> 50, K26: 30000 -> 25000, L26: 32000 -> 26000, M26: 31000 -> 25600, P26: 2385 -> 1969
$ws.Range("D26").Value = 44313
$ws.Range("J26").Value = 50
$ws.Range("K26").Value = 25000
$ws.Range("L26").Value = 26000
$ws.Range("M26").Value = 25600
$ws.Range("P26").Value = 1969

# Row 27: D27: 44775 -> 44883, J27: 20 -> 60, K27: 12000 -> 14000, L27: 13000 -> 15000, M27: 12500 -> 14500, P27: 962 -> 1115
$ws.Range("D27").Value = 44883
$ws.Range("J27").Value = 60
$ws.Range("K27").Value = 14000
$ws.Range("L27").Value = 15000
$ws.Range("M27").Value = 14500
$ws.Range("P27").Value = 1115

# Row 28: D28: 44782 -> 44355, J28: 40 -> 60, K28: 13000 -> 18000, L28: 14000 -> 20000, M28: 13500 -> 19000, P28: 1038 -> 1462
$ws.Range("D28").Value = 44355
$ws.Range("J28").Value = 60
$ws.Range("K28").Value = 18000
$ws.Range("L28").Value = 20000
$ws.Range("M28").Value = 19000
$ws.Range("P28").Value = 1462

# Row 29: D29: 44523 -> 44488, K29: 15000 -> 16000, L29: 16000 -> 17000, M29: 15500 -> 16500, P29: 1192 -> 1269
$ws.Range("D29").Value = 44488
$ws.Range("K29").Value = 16000
$ws.Range("L29").Value = 17000
$ws.Range("M29").Value = 16500
$ws.Range("P29").Value = 1269

# Row 30: D30: 44462 -> 44377, J30: 60 -> 40
$ws.Range("D30").Value = 44377
$ws.Range("J30").Value = 40

# Row 31: D31: 44777 -> 44308, J31: 25 -> 50, K31: 13000 -> 26000, L31: 14000 -> 27000, M31: 13600 -> 26400, P31: 1046 -> 2031
$ws.Range("D31").Value = 44308
$ws.Range("J31").Value = 50
$ws.Range("K31").Value = 26000
$ws.Range("L31").Value = 27000
$ws.Range("M31").Value = 26400
$ws.Range("P31").Value = 2031

# Row 32: D32: 44769 -> 44719, K32: 14000 -> 13000, L32: 15000 -> 14000, M32: 14600 -> 13400, P32: 1123 -> 1031
$ws.Range("D32").Value = 44719
$ws.Range("K32").Value = 13000
$ws.Range("L32").Value = 14000
$ws.Range("M32").Value = 13400
$ws.Range("P32").Value = 1031

# Row 33: D33: 44316 -> 44462, J33: 50 -> 60, K33: 27000 -> 14000, L33: 28000 -> 15000, M33: 27400 -> 14500, P33: 2108 -> 1115
$ws.Range("D33").Value = 44462
$ws.Range("J33").Value = 60
$ws.Range("K33").Value = 14000
$ws.Range("L33").Value = 15000
$ws.Range("M33").Value = 14500
$ws.Range("P33").Value = 1115

# Row 34: D34: 44159 -> 44433, J34: 60 -> 100, K34: 30000 -> 13000, L34: 32000 -> 14000, M34: 31000 -> 13500, P34: 2385 -> 1038
$ws.Range("D34").Value = 44433
$ws.Range("J34").Value = 100
$ws.Range("K34").Value = 13000
$ws.Range("L34").Value = 14000
$ws.Range("M34").Value = 13500
$ws.Range("P34").Value = 1038

# Row 35: D35: 44838 -> 44777, J35: 40 -> 25, K35: 14000 -> 13000, L35: 15000 -> 14000, M35: 14500 -> 13600, P35: 1115 -> 1046
$ws.Range("D35").Value = 44777
$ws.Range("J35").Value = 25
$ws.Range("K35").Value = 13000
$ws.Range("L35").Value = 14000
$ws.Range("M35").Value = 13600
$ws.Range("P35").Value = 1046

# Row 36: D36: 44883 -> 44705, J36: 60 -> 50, K36: 14000 -> 10000, L36: 15000 -> 11000, M36: 14500 -> 10400, P36: 1115 -> 800
$ws.Range("D36").Value = 44705
$ws.Range("J36").Value = 50
$ws.Range("K36").Value = 10000
$ws.Range("L36").Value = 11000
$ws.Range("M36").Value = 10400
$ws.Range("P36").Value = 800

# Row 37: D37: 44755 -> 44383, J37: 40 -> 50, K37: 14000 -> 15000, L37: 15000 -> 16000, M37: 14500 -> 15400, P37: 1115 -> 1185
$ws.Range("D37").Value = 44383
$ws.Range("J37").Value = 50
$ws.Range("K37").Value = 15000
$ws.Range("L37").Value = 16000
$ws.Range("M37").Value = 15400
$ws.Range("P37").Value = 1185

# Row 38: D38: 44474 -> 44264, K38: 13000 -> 30000, L38: 14000 -> 32000, M38: 13500 -> 31000, P38: 1038 -> 2385
$ws.Range("D38").Value = 44264
$ws.Range("K38").Value = 30000
$ws.Range("L38").Value = 32000
$ws.Range("M38").Value = 31000
$ws.Range("P38").Value = 2385

# Row 39: D39: 44810 -> 44782, J39: 50 -> 40, K39: 11000 -> 13000, L39: 12000 -> 14000, M39: 11600 -> 13500, P39: 892 -> 1038
$ws.Range("D39").Value = 44782
$ws.Range("J39").Value = 40
$ws.Range("K39").Value = 13000
$ws.Range("L39").Value = 14000
$ws.Range("M39").Value = 13500
$ws.Range("P39").Value = 1038

# Row 40: D40: 44425 -> 44467, J40: 60 -> 100, K40: 14000 -> 13000, L40: 15000 -> 14000, M40: 14500 -> 13500, P40: 1115 -> 1038
$ws.Range("D40").Value = 44467
$ws.Range("J40").Value = 100
$ws.Range("K40").Value = 13000
$ws.Range("L40").Value = 14000
$ws.Range("M40").Value = 13500
$ws.Range("P40").Value = 1038

# Row 41: D41: 44467 -> 44510, J41: 100 -> 40, K41: 13000 -> 15000, L41: 14000 -> 16000, M41: 13500 -> 15500, P41: 1038 -> 1192
$ws.Range("D41").Value = 44510
$ws.Range("J41").Value = 40
$ws.Range("K41").Value = 15000
$ws.Range("L41").Value = 16000
$ws.Range("M41").Value = 15500
$ws.Range("P41").Value = 1192

# Row 42: D42: 44488 -> 44838, K42: 16000 -> 14000, L42: 17000 -> 15000, M42: 16500 -> 14500, P42: 1269 -> 1115
$ws.Range("D42").Value = 44838
$ws.Range("K42").Value = 14000
$ws.Range("L42").Value = 15000
$ws.Range("M42").Value = 14500
$ws.Range("P42").Value = 1115

# Row 43: D43: 44741 -> 44320, K43: 14000 -> 26000, L43: 15000 -> 28000, M43: 14400 -> 26800, P43: 1108 -> 2062
$ws.Range("D43").Value = 44320
$ws.Range("K43").Value = 26000
$ws.Range("L43").Value = 28000
$ws.Range("M43").Value = 26800
$ws.Range("P43").Value = 2062

# Row 44: D44: 44691 -> 44810, J44: 100 -> 50, K44: 12000 -> 11000, L44: 13000 -> 12000, M44: 12500 -> 11600, P44: 962 -> 892
$ws.Range("D44").Value = 44810
$ws.Range("J44").Value = 50
$ws.Range("K44").Value = 11000
$ws.Range("L44").Value = 12000
$ws.Range("M44").Value = 11600
$ws.Range("P44").Value = 892
